# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns,
# and fix the B/C/D/E ordering for ApeXProtocol / FraxShare / EnergySwap rows 45-47.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.919.65"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "2.300.99"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'306.55"
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("D6").Value = "'97.56"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "'0.512"
$ws.Range("E7").Value = "  -1.34%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.506"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").Value = "'35.76"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("D11").Value = "'0.0790"
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").Value = "'18.24"
$ws.Range("E12").Value = "  +1.71%  "
$ws.Range("D13").Value = "'0.119"
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").Value = "2.658.73"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "2.304.37"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "'0.786"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "42.847.08"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").Value = "'12.71"
$ws.Range("E19").Value = "  -3.79%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "'6.05"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").Value = "'67.85"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").Value = "'236.74"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("D24").Value = "'2.15"
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("D25").Value = "'2.47"
$ws.Range("E25").Value = "  +2.23%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").Value = "'25.42"
$ws.Range("E28").Value = "  +3.14%  "
$ws.Range("D29").Value = "'166.86"
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("D31").Value = "'9.06"
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").Value = "'33.21"
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").Value = "'4.83"
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("D35").Value = "'5.02"
$ws.Range("E35").Value = "  -2.29%  "
$ws.Range("E36").Value = "  -4.77%  "
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("D38").Value = "'0.0692"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("D42").Value = "'2.74"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "2.008.26"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'2.15"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'10.03"
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'18.01"
$ws.Range("E47").Value = "  +4.34%  "
$ws.Range("D48").Value = "'2.79"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("D49").Value = "'53.89"
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("D51").Value = "2.526.09"
$ws.Range("E51").Value = "  +0.10%  "
